# Insert a new weekly price record as row 285 in the "Apio" sheet,
# pushing the existing row 285 (and everything below it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 285 (and all following rows) down, leaving a fresh blank row 285.
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A285").Value = 10
$ws.Range("B285").Value = "Vega Modelo de Temuco"
$ws.Range("C285").Value = "La Araucanía"
$ws.Range("D285").Value2 = 44798
$ws.Range("E285").Value = 9
$ws.Range("F285").Value = 100112017
$ws.Range("G285").Value = "Apio"
$ws.Range("H285").Value = "Americana (o)"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 400
$ws.Range("K285").Value = 11000
$ws.Range("L285").Value = 13000
$ws.Range("M285").Value = 12000
$ws.Range("N285").Value = "$/docena de matas"
$ws.Range("O285").Value = "Provincia del Elquí"
$ws.Range("P285").Value = 2000
$ws.Range("Q285").Value = 6
$ws.Range("R285").Value = "Hortaliza"
